$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = 4.0
$ws.Range("L1").Value = 90.0
